# Applies the highlighting / run-merge changes described by the diff:
#  1) "...результаты тестов;" paragraph: merge its 4 runs into a single run
#     (via Find&Replace so the runs collapse into one) and highlight the
#     whole paragraph (incl. the paragraph mark) yellow.
#  2) "...выбор предмета, а затем выбор вида теста;" paragraph: highlight
#     the whole paragraph (incl. the paragraph mark) green.
#  3) "...создания или удаления теста по темам;" paragraph: highlight the
#     whole paragraph (incl. the paragraph mark) green.
#  4) "...видеть таймер во время ответов на вопросы;" paragraph: highlight
#     only the run text yellow, leaving the paragraph-mark run properties
#     (w:pPr/w:rPr) untouched.

$d = $word.ActiveDocument

# wdColorIndex constants
$wdYellow = 7
$wdBrightGreen = 4

function Find-ParagraphIndex($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text -like $needle) {
            return $i
        }
    }
    return -1
}

# --- 1) merge the split "результаты тестов" runs into one run, then
#        highlight the whole paragraph yellow -----------------------------
$idx1 = Find-ParagraphIndex $d "*результаты тестов*"
if ($idx1 -lt 0) { throw "Could not locate 'результаты тестов' paragraph" }
$p1 = $d.Paragraphs.Item($idx1)
$p1.Range.Find.Execute(
    "предоставить преподавателю и студенту результаты тестов;", $true, $false,
    $false, $false, $false, $true, 1, $false,
    "предоставить преподавателю и студенту результаты тестов;", 2) | Out-Null

$idx1 = Find-ParagraphIndex $d "*результаты тестов*"
$p1 = $d.Paragraphs.Item($idx1)
$p1.Range.Font.HighlightColorIndex = $wdYellow

# --- 2) "выбор предмета, а затем выбор вида теста;" -> green, whole para --
$idx2 = Find-ParagraphIndex $d "*выбор предмета, а затем выбор*"
if ($idx2 -lt 0) { throw "Could not locate 'выбор предмета' paragraph" }
$p2 = $d.Paragraphs.Item($idx2)
$p2.Range.Font.HighlightColorIndex = $wdBrightGreen

# --- 3) "создания или удаления теста по темам;" -> green, whole para -----
$idx3 = Find-ParagraphIndex $d "*создания или удаления*"
if ($idx3 -lt 0) { throw "Could not locate 'создания или удаления' paragraph" }
$p3 = $d.Paragraphs.Item($idx3)
$p3.Range.Font.HighlightColorIndex = $wdBrightGreen

# --- 4) "видеть таймер во время ответов на вопросы;" -> yellow, runs only,
#        paragraph mark (w:pPr/w:rPr) must stay unhighlighted -------------
$idx4 = Find-ParagraphIndex $d "*видеть таймер*"
if ($idx4 -lt 0) { throw "Could not locate 'видеть таймер' paragraph" }
$p4 = $d.Paragraphs.Item($idx4)
$r4 = $p4.Range
$textOnly4 = $d.Range($r4.Start, $r4.End - 1)
$textOnly4.Font.HighlightColorIndex = $wdYellow
